$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 233-244: only the cells that actually changed ---

# Row 233
$ws.Range("D233").Value = 44931
$ws.Range("K233").Value = "Bing"
$ws.Range("M233").Value = 500
$ws.Range("N233").Value = 7000
$ws.Range("O233").Value = 8000
$ws.Range("P233").Value = 7500
$ws.Range("R233").Value = "Provincia de Curicó"
$ws.Range("S233").Value = 750

# Row 234
$ws.Range("D234").Value = 44931
$ws.Range("K234").Value = "Bing"
$ws.Range("N234").Value = 5000
$ws.Range("O234").Value = 6000
$ws.Range("P234").Value = 5500
$ws.Range("R234").Value = "Provincia de Curicó"
$ws.Range("S234").Value = 550

# Row 235
$ws.Range("D235").Value = 44931
$ws.Range("K235").Value = "Lapins"
$ws.Range("N235").Value = 11000
$ws.Range("O235").Value = 12000
$ws.Range("P235").Value = 11500
$ws.Range("S235").Value = 767

# Row 236
$ws.Range("D236").Value = 44931
$ws.Range("K236").Value = "Lapins"
$ws.Range("L236").Value = "Segunda"
$ws.Range("M236").Value = 400
$ws.Range("N236").Value = 9000
$ws.Range("O236").Value = 10000
$ws.Range("P236").Value = 9500
$ws.Range("S236").Value = 633

# Row 237
$ws.Range("K237").Value = "Lapins"
$ws.Range("L237").Value = "Primera"
$ws.Range("M237").Value = 340
$ws.Range("N237").Value = 11000
$ws.Range("O237").Value = 12000
$ws.Range("P237").Value = 11500
$ws.Range("S237").Value = 1150

# Row 238
$ws.Range("K238").Value = "Lapins"
$ws.Range("L238").Value = "Segunda"
$ws.Range("M238").Value = 400
$ws.Range("N238").Value = 9000
$ws.Range("O238").Value = 10000
$ws.Range("P238").Value = 9500
$ws.Range("S238").Value = 950

# Row 239
$ws.Range("K239").Value = "Rainier"
$ws.Range("L239").Value = "Primera"
$ws.Range("M239").Value = 400
$ws.Range("N239").Value = 20000
$ws.Range("O239").Value = 21000
$ws.Range("P239").Value = 20500
$ws.Range("Q239").Value = "`$/caja 15 kilos"
$ws.Range("R239").Value = "Provincia de Curicó"
$ws.Range("S239").Value = 1367
$ws.Range("T239").Value = 15

# Row 240
$ws.Range("D240").Value = 44545
$ws.Range("K240").Value = "Royal Dawn"
$ws.Range("M240").Value = 300
$ws.Range("N240").Value = 14000
$ws.Range("O240").Value = 15000
$ws.Range("P240").Value = 14500
$ws.Range("Q240").Value = "`$/caja 15 kilos"
$ws.Range("S240").Value = 967
$ws.Range("T240").Value = 15

# Row 241
$ws.Range("D241").Value = 44545
$ws.Range("K241").Value = "Santina"
$ws.Range("L241").Value = "Especial"
$ws.Range("M241").Value = 360
$ws.Range("N241").Value = 13000
$ws.Range("O241").Value = 14000
$ws.Range("P241").Value = 13500
$ws.Range("R241").Value = "Región de O'Higgins"
$ws.Range("S241").Value = 1350

# Row 242
$ws.Range("D242").Value = 44545
$ws.Range("M242").Value = 340
$ws.Range("N242").Value = 11000
$ws.Range("O242").Value = 12000
$ws.Range("P242").Value = 11500
$ws.Range("R242").Value = "Región de O'Higgins"
$ws.Range("S242").Value = 1150

# Row 243
$ws.Range("D243").Value = 44545
$ws.Range("M243").Value = 300
$ws.Range("N243").Value = 9000
$ws.Range("O243").Value = 10000
$ws.Range("P243").Value = 9500
$ws.Range("R243").Value = "Región de O'Higgins"
$ws.Range("S243").Value = 950

# Row 244
$ws.Range("K244").Value = "Lapins"
$ws.Range("N244").Value = 8500
$ws.Range("O244").Value = 9000
$ws.Range("P244").Value = 8750
$ws.Range("S244").Value = 875

# --- Append new rows 245-248 (full rows, matching the established column layout) ---

# Row 245
$ws.Range("A245").Value = 2
$ws.Range("B245").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44580
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100103
$ws.Range("H245").Value = "Frutos de hueso (carozo)"
$ws.Range("I245").Value = 100103001
$ws.Range("J245").Value = "Cereza"
$ws.Range("K245").Value = "Lapins"
$ws.Range("L245").Value = "Segunda"
$ws.Range("M245").Value = 400
$ws.Range("N245").Value = 6500
$ws.Range("O245").Value = 7000
$ws.Range("P245").Value = 6750
$ws.Range("Q245").Value = "`$/bandeja 10 kilos"
$ws.Range("R245").Value = "Provincia de Curicó"
$ws.Range("S245").Value = 675
$ws.Range("T245").Value = 10
$ws.Range("D245").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 246
$ws.Range("A246").Value = 2
$ws.Range("B246").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = 44580
$ws.Range("E246").Value = 4
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100103
$ws.Range("H246").Value = "Frutos de hueso (carozo)"
$ws.Range("I246").Value = 100103001
$ws.Range("J246").Value = "Cereza"
$ws.Range("K246").Value = "Santina"
$ws.Range("L246").Value = "Primera"
$ws.Range("M246").Value = 500
$ws.Range("N246").Value = 8500
$ws.Range("O246").Value = 9000
$ws.Range("P246").Value = 8750
$ws.Range("Q246").Value = "`$/bandeja 10 kilos"
$ws.Range("R246").Value = "Provincia de Curicó"
$ws.Range("S246").Value = 875
$ws.Range("T246").Value = 10
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 247
$ws.Range("A247").Value = 2
$ws.Range("B247").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44580
$ws.Range("E247").Value = 4
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100103
$ws.Range("H247").Value = "Frutos de hueso (carozo)"
$ws.Range("I247").Value = 100103001
$ws.Range("J247").Value = "Cereza"
$ws.Range("K247").Value = "Santina"
$ws.Range("L247").Value = "Segunda"
$ws.Range("M247").Value = 400
$ws.Range("N247").Value = 6500
$ws.Range("O247").Value = 7000
$ws.Range("P247").Value = 6750
$ws.Range("Q247").Value = "`$/bandeja 10 kilos"
$ws.Range("R247").Value = "Provincia de Curicó"
$ws.Range("S247").Value = 675
$ws.Range("T247").Value = 10
$ws.Range("D247").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 248
$ws.Range("A248").Value = 2
$ws.Range("B248").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C248").Value = "Coquimbo"
$ws.Range("D248").Value = 44580
$ws.Range("E248").Value = 4
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100103
$ws.Range("H248").Value = "Frutos de hueso (carozo)"
$ws.Range("I248").Value = 100103001
$ws.Range("J248").Value = "Cereza"
$ws.Range("K248").Value = "Sweet Heart"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 500
$ws.Range("N248").Value = 9500
$ws.Range("O248").Value = 10000
$ws.Range("P248").Value = 9750
$ws.Range("Q248").Value = "`$/bandeja 10 kilos"
$ws.Range("R248").Value = "Provincia de Curicó"
$ws.Range("S248").Value = 975
$ws.Range("T248").Value = 10
$ws.Range("D248").NumberFormat = "YYYY-MM-DD HH:MM:SS"
